$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 797.4
$ws.Range("I41").Value = 1001
$ws.Range("J41").Value = 746.5
$ws.Range("K41").Value = 1001
$ws.Range("L41").Value = 746.5
$ws.Range("M41").Value = -561
$ws.Range("N41").Value = -1626.5
$ws.Range("H74").Value = 3486.8333
$ws.Range("I74").Value = 3486.8333
$ws.Range("K74").Value = 3486.8333
$ws.Range("M74").Value = -2550.8333
$ws.Range("H77").Value = 3486.8333
$ws.Range("I77").Value = 3486.8333
$ws.Range("K77").Value = 17434.1665
$ws.Range("M77").Value = -12754.1665
$ws.Range("H100").Value = 2791.9285
$ws.Range("I100").Value = 2682.3333
$ws.Range("K100").Value = 2682.3333
$ws.Range("M100").Value = -2141.3333
$ws.Range("H111").Value = 1146
$ws.Range("I111").Value = 837
$ws.Range("K111").Value = 2511
$ws.Range("M111").Value = 556
$ws.Range("H113").Value = 5005
$ws.Range("I113").Value = 5005
$ws.Range("K113").Value = 5005
$ws.Range("M113").Value = -1751

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1311.2307
$ws.Range("I2").Value = 960.7778
$ws.Range("K2").Value = 960.7778
$ws.Range("M2").Value = -847.7778
$ws.Range("H6").Value = 1017500
$ws.Range("I6").Value = 1017500
$ws.Range("K6").Value = 1017500
$ws.Range("M6").Value = -1017327
$ws.Range("H32").Value = 4093.2
$ws.Range("I32").Value = 3139.1785
$ws.Range("K32").Value = 3139.1785
$ws.Range("M32").Value = -2852.1785
$ws.Range("H37").Value = 6356
$ws.Range("I37").Value = 6356
$ws.Range("K37").Value = 6356
$ws.Range("M37").Value = -6083
$ws.Range("H116").Value = 1311.2307
$ws.Range("I116").Value = 960.7778
$ws.Range("K116").Value = 960.7778
$ws.Range("M116").Value = 1333.2222
$ws.Range("H131").Value = 64856.25
$ws.Range("J131").Value = 64856.25
$ws.Range("L131").Value = 64856.25
$ws.Range("N131").Value = -74936.25
$ws.Range("H132").Value = 5249.96
$ws.Range("I132").Value = 5488.2856
$ws.Range("J132").Value = 3998.75
$ws.Range("K132").Value = 16464.8568
$ws.Range("L132").Value = 11996.25
$ws.Range("M132").Value = -13934.8568
$ws.Range("N132").Value = -17056.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1311.2307
$ws.Range("I3").Value = 960.7778
$ws.Range("K3").Value = 960.7778
$ws.Range("M3").Value = -846.7778
$ws.Range("H7").Value = 2800
$ws.Range("J7").Value = 2800
$ws.Range("L7").Value = 2800
$ws.Range("N7").Value = -3026
$ws.Range("H20").Value = 3316.6667
$ws.Range("I20").Value = 2756
$ws.Range("J20").Value = 3765.2
$ws.Range("K20").Value = 2756
$ws.Range("L20").Value = 3765.2
$ws.Range("M20").Value = -2509
$ws.Range("N20").Value = -4259.2
$ws.Range("H96").Value = 17999
$ws.Range("I96").Value = 17999
$ws.Range("K96").Value = 17999
$ws.Range("M96").Value = -15253
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H134").Value = 5877.077
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1926.3
$ws.Range("I31").Value = 1908.5
$ws.Range("K31").Value = 1908.5
$ws.Range("M31").Value = -1613.5
$ws.Range("H33").Value = 1437.8334
$ws.Range("I33").Value = 1619.6
$ws.Range("J33").Value = 529
$ws.Range("K33").Value = 1619.6
$ws.Range("L33").Value = 529
$ws.Range("M33").Value = -1240.6
$ws.Range("N33").Value = -1287
$ws.Range("H34").Value = 1926.3
$ws.Range("I34").Value = 1908.5
$ws.Range("K34").Value = 1908.5
$ws.Range("M34").Value = -1706.5
$ws.Range("H86").Value = 8081.5
$ws.Range("I86").Value = 7489
$ws.Range("K86").Value = 7489
$ws.Range("M86").Value = -6366
$ws.Range("H89").Value = 8081.5
$ws.Range("I89").Value = 7489
$ws.Range("K89").Value = 37445
$ws.Range("M89").Value = -31829
$ws.Range("H105").Value = 2294
$ws.Range("I105").Value = 1989.1111
$ws.Range("J105").Value = 3666
$ws.Range("K105").Value = 1989.1111
$ws.Range("L105").Value = 3666
$ws.Range("M105").Value = -242.1111000000001
$ws.Range("N105").Value = -7160
$ws.Range("H107").Value = 1427.2727
$ws.Range("I107").Value = 608.8570999999999
$ws.Range("J107").Value = 2859.5
$ws.Range("K107").Value = 608.8570999999999
$ws.Range("L107").Value = 2859.5
$ws.Range("M107").Value = 1311.1429
$ws.Range("N107").Value = -6699.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 975
$ws.Range("I25").Value = 900
$ws.Range("K25").Value = 2700
$ws.Range("M25").Value = -2531
$ws.Range("H30").Value = 975
$ws.Range("I30").Value = 900
$ws.Range("K30").Value = 2700
$ws.Range("M30").Value = -2598
$ws.Range("H40").Value = 50
$ws.Range("I40").Value = 50
$ws.Range("K40").Value = 200
$ws.Range("M40").Value = -131
$ws.Range("H129").Value = 2813.6667
$ws.Range("J129").Value = 2813.6667
$ws.Range("L129").Value = 8441.000100000001
$ws.Range("N129").Value = -18441.0001
$ws.Range("H133").Value = 17103.5
$ws.Range("J133").Value = 24333.334
$ws.Range("L133").Value = 73000.00199999999
$ws.Range("N133").Value = -83120.00199999999
$ws.Range("H134").Value = 1850.8
$ws.Range("I134").Value = 1063.5
$ws.Range("K134").Value = 3190.5
$ws.Range("M134").Value = 1879.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 81.72727
$ws.Range("I2").Value = 87.90000000000001
$ws.Range("K2").Value = 87.90000000000001
$ws.Range("M2").Value = 25.09999999999999
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H102").Value = 2088.6667
$ws.Range("I102").Value = 2099.75
$ws.Range("K102").Value = 2099.75
$ws.Range("M102").Value = -477.75
$ws.Range("H132").Value = 4957.5835
$ws.Range("I132").Value = 4056.125
$ws.Range("K132").Value = 12168.375
$ws.Range("M132").Value = -9638.375
$ws.Range("H136").Value = 15769.75
$ws.Range("J136").Value = 15769.75
$ws.Range("L136").Value = 47309.25
$ws.Range("N136").Value = -52409.25
$ws.Range("M58").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 757
$ws.Range("I16").Value = 579.4286
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 579.4286
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -409.4286
$ws.Range("N16").Value = -2340
$ws.Range("H40").Value = 3229.75
$ws.Range("I40").Value = 3229.75
$ws.Range("K40").Value = 3229.75
$ws.Range("M40").Value = -3093.75
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H61").Value = 1861.5
$ws.Range("I61").Value = 1861.5
$ws.Range("K61").Value = 1861.5
$ws.Range("M61").Value = -1659.5
$ws.Range("H68").Value = 2500
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1751
$ws.Range("H71").Value = 2500
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -8756
$ws.Range("H93").Value = 3750
$ws.Range("I93").Value = 3750
$ws.Range("K93").Value = 3750
$ws.Range("M93").Value = -2502
$ws.Range("H113").Value = 1861.5
$ws.Range("I113").Value = 1861.5
$ws.Range("K113").Value = 1861.5
$ws.Range("M113").Value = 308.5
$ws.Range("M57").ClearContents()
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 57497.5
$ws.Range("J130").Value = 57497.5
$ws.Range("L130").Value = 57497.5
$ws.Range("N130").Value = -67537.5
